$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Delete-Text($old) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
}

# --- Title / author / email ---
Replace-Text "Cosmic Enlightenment: Hubble's Visionary Journey" "The Art of Numbers: A Journey into the Realm of Mathematics"
Replace-Text "Amelia Reynolds" "Alexia Hartwell"
Replace-Text "areynolds@stellarobservatory" "alexiahartwell@validmail"

# --- Body paragraph 1, block 1 ---
Replace-Text "With every click of the camera shutter, Edwin Hubble carved an ineffaceable legacy, reshaping our understanding of the universe we inhabit" "In the vast expanse of human knowledge, Mathematics reigns supreme as the language of science, reason, and order"
Replace-Text " His pioneering work, fueled by the insatiable quest for knowledge enshrined within his gaze, transformed the 20th-century scientific landscape" " It is a discipline that has captivated minds for millennia, unraveling the intricate patterns that weave the fabric of our universe"
Replace-Text " Hubble's telescope, a colossal eye reaching beyond the confines of our earthly abode, bore witness to the boundless panorama of the cosmos, revealing secrets hitherto shrouded in the cloak of obscurity" " From the awe-inspiring cosmos to the intricate structures of living organisms, Mathematics provides a powerful lens through which we can understand and interpret the world around us"

# Remove the next three now-orphaned sentences (each still prefixed with its own leading space, the previous sentence's trailing period stays)
Delete-Text " The observable universe expanded under his scrutiny, transforming static cosmic maps into dynamic tapestries of evolution and motion. Hubble's observations laid bare the existence of galaxies beyond our Milky Way, shattering the long-held belief in our celestial solitude. His discoveries illuminated the vastness of the universe, kindling awe and wonder within the hearts of humanity."

# --- Body paragraph 1, block 2 ---
Replace-Text "Propelled by an unwavering spirit of curiosity, Hubble embarked on a tireless quest to unravel the enigmas of the cosmos" "Delving into the realm of Mathematics is akin to embarking on an enthralling adventure, where the exploration of numbers, patterns, and relationships unlocks hidden mysteries"
Replace-Text " His meticulous observations fueled a surge of scientific exploration and discovery, propelling astronomers to the forefront of our collective understanding of the universe" " It is a journey that requires curiosity, creativity, and a willingness to embrace the challenges that lie ahead"
Replace-Text " With each celestial snapshot, Hubble unveiled new realms of knowledge, challenging prevailing theories and revolutionizing our comprehension of cosmic phenomena" " As we navigate this intellectual landscape, we uncover profound truths about the nature of reality, the limits of our understanding, and the immense possibilities that lie within the realm of human thought"

Delete-Text " His groundbreaking studies elucidated the existence of red-shifted galaxies, providing irrefutable evidence for the expansion of the universe and paving the way for the revolutionary theory of the Big Bang. Hubble's contributions extended far beyond his lifetime, influencing generations of astronomers and fueling a legacy of tireless exploration."

# --- Body paragraph 1, block 3 ---
Replace-Text "Hubble's unwavering commitment to unveiling the mysteries of the universe serves as a testament to the transformative power of human curiosity" "Mathematics is not merely a collection of abstract concepts; it is a vibrant and dynamic discipline that finds its applications in every corner of our lives"
Replace-Text " His pioneering spirit and relentless pursuit of knowledge continue to inspire scientists and stargazers alike, igniting a flame of wonder and exploration that transcends generations" " Whether it be the construction of towering skyscrapers, the intricate designs of electronic circuits, or the complex algorithms that drive our digital age, Mathematics plays an indispensable role in shaping our world"
Replace-Text " Through his extraordinary endeavors, Hubble not only redefined our comprehension " " It is a tool that empowers us to solve problems, make predictions, and explore the boundless frontiers of human knowledge"
Replace-Text "of the cosmos but also ignited an eternal quest for understanding our place within this vast and awe-inspiring universe" "numbers that describe the intricate relationships between abstract quantities, numbers are ubiquitous in our daily lives"

# Insert the new "Body:" / "Paragraph 1/2/3" sections right after the sentence that now ends "...human knowledge."
$rng = $d.Content
$rng.Find.Execute("It is a tool that empowers us to solve problems, make predictions, and explore the boundless frontiers of human knowledge") | Out-Null
$rng.Collapse(0)
$newline = [char]11
$insert1 = ".$newline${newline}Body:$newline${newline}Paragraph 1:$newline" + `
  "At the heart of Mathematics lies the concept of numbers, those fundamental building blocks that form the foundation of all mathematical thought. " + `
  "From the simple counting numbers that we use to enumerate objects to the complex "
$rng.InsertAfter($insert1)

# The text that follows (after the lastRenderedPageBreak run) becomes the continuation of paragraph 1, plus paragraphs 2 and 3
$rng2 = $d.Content
$rng2.Find.Execute("numbers that describe the intricate relationships between abstract quantities, numbers are ubiquitous in our daily lives") | Out-Null
$rng2.Collapse(0)
$insert2 = ". We use them to measure, compare, calculate, and understand the world around us.$newline${newline}Paragraph 2:$newline" + `
  "Beyond numbers, Mathematics encompasses a vast array of concepts, theorems, and principles that govern the interactions between quantities and shapes. " + `
  "Algebra, with its abstract symbols and equations, provides a powerful framework for solving complex problems. " + `
  "Geometry, with its focus on shapes and spatial relationships, helps us visualize and understand the world in three dimensions. " + `
  "Calculus, with its intricate techniques for analyzing change, unlocks the secrets of motion and growth.$newline${newline}Paragraph 3:$newline" + `
  "The true beauty of Mathematics lies in its universality. " + `
  "It transcends cultural, linguistic, and geographical boundaries, uniting people from all corners of the globe in a shared pursuit of knowledge. " + `
  "Mathematical concepts and principles are the same for everyone, regardless of their background or beliefs. " + `
  "This universality makes Mathematics a powerful tool for communication and collaboration, enabling scientists, engineers, and mathematicians from around the world to work together to solve global challenges"
$rng2.InsertAfter($insert2)

# --- Summary paragraph ---
Replace-Text "Edwin Hubble's visionary journey through the cosmos, armed with his pioneering telescope, transformed our comprehension of the universe" "Mathematics is a fascinating and multifaceted discipline that offers profound insights into the nature of reality"
Replace-Text " His meticulous observations revealed the existence of galaxies beyond our own, propelling us into an era of profound scientific discovery" " Through the study of numbers, patterns, and relationships, we gain a deeper understanding of the world around us and the boundless possibilities that lie within the realm of human thought"
Replace-Text " Hubble's groundbreaking studies elucidated the expansion of the universe, laying the foundation for the Big Bang theory and revolutionizing our understanding of cosmic evolution" " Its universality makes it a powerful tool for communication and collaboration, enabling people from all over the world to work together to solve complex problems"
Replace-Text " His unwavering spirit of curiosity and relentless pursuit of knowledge continue to inspire generations, fueling a legacy of tireless exploration that transcends time and space" " Mathematics is not just a subject that we learn in school; it is a way of thinking, a tool for exploration, and a window into the intricate workings of the universe"

Delete-Text ". Hubble's name remains synonymous with innovation, ingenuity, and the transformative power of human curiosity"

# --- Add a new empty paragraph at the end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
